# Applies the "ICAS: End of year run 2024" update to the model-coefficient
# lookup tables on sheets Q100-Q600, P100-P600 and IK106, refreshing the fitted
# coefficients, the effective starting dates and the "gewijzigd" timestamps.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("IK106")
$ws.Range("E2").Value = [double]"45096.51158341435"
$ws.Range("E3").Value = [double]"45096.51158341435"
$ws.Range("E4").Value = [double]"45096.51158357639"
$ws.Range("E5").Value = [double]"45096.51158341435"

$ws = $wb.Worksheets.Item("Q100")
$ws.Range("B2").Value = [double]"41030.97916666666"
$ws.Range("C2").Value = [double]"-9.812742932787027e-06"
$ws.Range("D2").Value = [double]"-1.2e-09"
$ws.Range("E2").Value = [double]"45657.5382174537"
$ws.Range("C3").Value = [double]"-1.219416793278703e-05"
$ws.Range("D3").Value = [double]"-1.2e-09"
$ws.Range("E3").Value = [double]"45657.53821756945"
$ws.Range("C4").Value = [double]"-1.249929159773437e-05"
$ws.Range("D4").Value = [double]"-1.2e-09"
$ws.Range("E4").Value = [double]"45657.5382174537"
$ws.Range("C5").Value = [double]"-1.146513897339869e-05"
$ws.Range("D5").Value = [double]"-1.2e-09"
$ws.Range("E5").Value = [double]"45657.5382174537"

$ws = $wb.Worksheets.Item("Q200")
$ws.Range("C2").Value = [double]"-1.482359601319953e-05"
$ws.Range("D2").Value = [double]"-4.99e-10"
$ws.Range("E2").Value = [double]"45657.53840547454"
$ws.Range("C3").Value = [double]"-1.548383080773783e-05"
$ws.Range("D3").Value = [double]"-4.99e-10"
$ws.Range("E3").Value = [double]"45657.53840547454"
$ws.Range("C4").Value = [double]"-1.546812433875638e-05"
$ws.Range("D4").Value = [double]"-4.99e-10"
$ws.Range("E4").Value = [double]"45657.53840547454"

$ws = $wb.Worksheets.Item("Q300")
$ws.Range("C2").Value = [double]"-7.743839100585073e-06"
$ws.Range("D2").Value = [double]"-1.32e-09"
$ws.Range("E2").Value = [double]"45657.53867140046"
$ws.Range("C3").Value = [double]"-8.890919100585073e-06"
$ws.Range("D3").Value = [double]"-1.32e-09"
$ws.Range("E3").Value = [double]"45657.53867146991"
$ws.Range("C4").Value = [double]"-1.032773910058507e-05"
$ws.Range("D4").Value = [double]"-1.32e-09"
$ws.Range("E4").Value = [double]"45657.53867146991"
$ws.Range("C5").Value = [double]"-7.081335203104816e-06"
$ws.Range("D5").Value = [double]"-1.32e-09"
$ws.Range("E5").Value = [double]"45657.53867140046"

$ws = $wb.Worksheets.Item("Q400")
$ws.Range("C2").Value = [double]"-9.253206191052623e-06"
$ws.Range("D2").Value = [double]"-4.76e-10"
$ws.Range("E2").Value = [double]"45657.53885164352"
$ws.Range("C3").Value = [double]"-9.384863941728345e-06"
$ws.Range("D3").Value = [double]"-4.76e-10"
$ws.Range("E3").Value = [double]"45657.53885164352"
$ws.Range("C4").Value = [double]"-8.864197589164204e-06"
$ws.Range("D4").Value = [double]"-4.76e-10"
$ws.Range("E4").Value = [double]"45657.53885164352"
$ws.Range("C5").Value = [double]"-8.253573208534444e-06"
$ws.Range("D5").Value = [double]"-4.76e-10"
$ws.Range("E5").Value = [double]"45657.53885164352"

$ws = $wb.Worksheets.Item("Q500")
$ws.Range("C2").Value = [double]"-1.567735900203139e-05"
$ws.Range("D2").Value = [double]"-8.12e-10"
$ws.Range("E2").Value = [double]"45657.53912391204"
$ws.Range("C3").Value = [double]"-1.670981700203139e-05"
$ws.Range("D3").Value = [double]"-8.12e-10"
$ws.Range("E3").Value = [double]"45657.53912417824"
$ws.Range("C4").Value = [double]"-1.759083700203139e-05"
$ws.Range("D4").Value = [double]"-8.12e-10"
$ws.Range("E4").Value = [double]"45657.53912417824"
$ws.Range("C5").Value = [double]"-1.575197772123149e-05"
$ws.Range("D5").Value = [double]"-8.12e-10"
$ws.Range("E5").Value = [double]"45657.53912391204"
$ws.Range("C6").Value = [double]"-8.814434289001204e-06"
$ws.Range("D6").Value = [double]"-8.12e-10"
$ws.Range("E6").Value = [double]"45657.53912391204"

$ws = $wb.Worksheets.Item("Q600")
$ws.Range("C2").Value = [double]"-1.735425541454964e-05"
$ws.Range("E2").Value = [double]"45657.53933913194"
$ws.Range("C3").Value = [double]"-1.918275541454964e-05"
$ws.Range("E3").Value = [double]"45657.53933918982"
$ws.Range("C4").Value = [double]"-1.934179400942103e-05"
$ws.Range("E4").Value = [double]"45657.53933913194"
$ws.Range("C5").Value = [double]"-2.002349721098186e-05"
$ws.Range("E5").Value = [double]"45657.53933913194"

$ws = $wb.Worksheets.Item("P100")
$ws.Range("B2").Value = [double]"41030.97916666666"
$ws.Range("C2").Value = [double]"-1.141922310494986e-05"
$ws.Range("D2").Value = [double]"-4.7e-10"
$ws.Range("E2").Value = [double]"45657.53954047454"
$ws.Range("C3").Value = [double]"-9.49886197188645e-06"
$ws.Range("D3").Value = [double]"-4.7e-10"
$ws.Range("E3").Value = [double]"45657.53954047454"
$ws.Range("C4").Value = [double]"-1.006803197188645e-05"
$ws.Range("D4").Value = [double]"-4.7e-10"
$ws.Range("E4").Value = [double]"45657.53954060185"
$ws.Range("C5").Value = [double]"-9.514804369884435e-06"
$ws.Range("D5").Value = [double]"-4.7e-10"
$ws.Range("E5").Value = [double]"45657.53954047454"

$ws = $wb.Worksheets.Item("P200")
$ws.Range("C2").Value = [double]"-1.230349195937308e-05"
$ws.Range("D2").Value = [double]"-7.7e-10"
$ws.Range("E2").Value = [double]"45657.53980603009"
$ws.Range("C3").Value = [double]"-1.188855053699145e-05"
$ws.Range("D3").Value = [double]"-7.7e-10"
$ws.Range("E3").Value = [double]"45657.53980603009"
$ws.Range("C4").Value = [double]"-1.246797553699145e-05"
$ws.Range("D4").Value = [double]"-7.7e-10"
$ws.Range("E4").Value = [double]"45657.53980618055"
$ws.Range("C5").Value = [double]"-1.179156717366791e-05"
$ws.Range("D5").Value = [double]"-7.7e-10"
$ws.Range("E5").Value = [double]"45657.53980603009"
$ws.Range("C6").Value = [double]"-1.098257305276865e-05"
$ws.Range("D6").Value = [double]"-7.7e-10"
$ws.Range("E6").Value = [double]"45657.53980603009"

$ws = $wb.Worksheets.Item("P300")
$ws.Range("B2").Value = [double]"41035.97916666666"
$ws.Range("C2").Value = [double]"-1.234103346351092e-05"
$ws.Range("D2").Value = [double]"-1.23e-09"
$ws.Range("E2").Value = [double]"45657.54015164352"
$ws.Range("C3").Value = [double]"-1.313563908851092e-05"
$ws.Range("D3").Value = [double]"-1.23e-09"
$ws.Range("E3").Value = [double]"45657.54015172453"
$ws.Range("C4").Value = [double]"-1.391914908851092e-05"
$ws.Range("D4").Value = [double]"-1.23e-09"
$ws.Range("E4").Value = [double]"45657.54015172453"
$ws.Range("C5").Value = [double]"-1.482750408851092e-05"
$ws.Range("D5").Value = [double]"-1.23e-09"
$ws.Range("E5").Value = [double]"45657.54015172453"
$ws.Range("C6").Value = [double]"-1.333925029812126e-05"
$ws.Range("D6").Value = [double]"-1.23e-09"
$ws.Range("E6").Value = [double]"45657.54015164352"
$ws.Range("C7").Value = [double]"-1.354918647600766e-05"
$ws.Range("D7").Value = [double]"-1.23e-09"
$ws.Range("E7").Value = [double]"45657.54015164352"
$ws.Range("C8").Value = [double]"-1.226796616201782e-05"
$ws.Range("D8").Value = [double]"-1.23e-09"
$ws.Range("E8").Value = [double]"45657.54015164352"

$ws = $wb.Worksheets.Item("P400")
$ws.Range("C2").Value = [double]"-6.111882445032339e-06"
$ws.Range("E2").Value = [double]"45657.5403791088"
$ws.Range("C3").Value = [double]"-5.616123389944509e-06"
$ws.Range("E3").Value = [double]"45657.5403791088"
$ws.Range("C4").Value = [double]"-7.125323389944509e-06"
$ws.Range("E4").Value = [double]"45657.54037918981"
$ws.Range("C5").Value = [double]"-5.827565334607773e-06"
$ws.Range("E5").Value = [double]"45657.5403791088"
$ws.Range("C6").Value = [double]"-5.217831045153544e-06"
$ws.Range("E6").Value = [double]"45657.5403791088"

$ws = $wb.Worksheets.Item("P500")
$ws.Range("C2").Value = [double]"-6.018697437020036e-06"
$ws.Range("D2").Value = [double]"-1.37e-09"
$ws.Range("E2").Value = [double]"45657.54068890047"
$ws.Range("C3").Value = [double]"-7.415412437020036e-06"
$ws.Range("D3").Value = [double]"-1.37e-09"
$ws.Range("E3").Value = [double]"45657.54068903935"
$ws.Range("C4").Value = [double]"-9.328617437020036e-06"
$ws.Range("D4").Value = [double]"-1.37e-09"
$ws.Range("E4").Value = [double]"45657.54068903935"
$ws.Range("C5").Value = [double]"-5.83149117296834e-06"
$ws.Range("D5").Value = [double]"-1.37e-09"
$ws.Range("E5").Value = [double]"45657.54068890047"

$ws = $wb.Worksheets.Item("P600")
$ws.Range("B2").Value = [double]"41030.97916666666"
$ws.Range("C2").Value = [double]"-8.00106841502271e-06"
$ws.Range("E2").Value = [double]"45657.54090913876"
$ws.Range("C3").Value = [double]"-8.5780375932064e-06"
$ws.Range("E3").Value = [double]"45657.54090913876"
$ws.Range("C4").Value = [double]"-1.197094368682955e-05"
$ws.Range("E4").Value = [double]"45657.54090913876"
$ws.Range("C5").Value = [double]"-1.016169640946589e-05"
$ws.Range("E5").Value = [double]"45657.54090913876"
